$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = "'66.538.95"
$ws.Cells.Item(2,4).ClearFormats()
$ws.Cells.Item(2,5).Value2 = '  -4.97%  '

$ws.Cells.Item(3,4).Value2 = "'3.223.45"
$ws.Cells.Item(3,4).ClearFormats()
$ws.Cells.Item(3,5).Value2 = '  -7.89%  '

$ws.Cells.Item(4,4).Value2 = "'1.01"
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).Value2 = '  +0.86%  '

$ws.Cells.Item(5,4).Value2 = "'585.65"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value2 = '  -2.76%  '

$ws.Cells.Item(6,4).Value2 = "'152.14"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value2 = '  -12.60%  '

$ws.Cells.Item(7,5).Value2 = '  +0.34%  '

$ws.Cells.Item(8,4).Value2 = "'3.207.79"
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value2 = '  -8.19%  '

$ws.Cells.Item(9,4).Value2 = "'0.542"
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value2 = '  -10.92%  '

$ws.Cells.Item(10,4).Value2 = "'0.173"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value2 = '  -9.79%  '

$ws.Cells.Item(11,4).Value2 = "'5.75"
$ws.Cells.Item(11,4).ClearFormats()
$ws.Cells.Item(11,5).Value2 = '  -20.57%  '

$ws.Cells.Item(12,4).Value2 = "'0.493"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value2 = '  -14.82%  '

$ws.Cells.Item(13,4).Value2 = "'38.73"
$ws.Cells.Item(13,4).ClearFormats()
$ws.Cells.Item(13,5).Value2 = '  -15.79%  '

$ws.Cells.Item(14,4).Value2 = "'0.0000240"
$ws.Cells.Item(14,4).ClearFormats()
$ws.Cells.Item(14,5).Value2 = '  -12.23%  '

$ws.Cells.Item(15,4).Value2 = "'3.837.24"
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).Value2 = '  -5.68%  '

$ws.Cells.Item(16,4).Value2 = "'67.171.11"
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value2 = '  -4.30%  '

$ws.Cells.Item(17,4).Value2 = "'3.306.23"
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value2 = '  -5.51%  '

$ws.Cells.Item(18,2).Value2 = 'TRON'
$ws.Cells.Item(18,3).Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(18,4).Value2 = "'0.115"
$ws.Cells.Item(18,4).ClearFormats()
$ws.Cells.Item(18,5).Value2 = '  -3.54%  '

$ws.Cells.Item(19,4).Value2 = "'534.12"
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value2 = '  -12.15%  '

$ws.Cells.Item(20,2).Value2 = 'Polkadot'
$ws.Cells.Item(20,3).Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(20,4).Value2 = "'7.11"
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value2 = '  -13.77%  '

$ws.Cells.Item(21,4).Value2 = "'14.98"
$ws.Cells.Item(21,4).ClearFormats()
$ws.Cells.Item(21,5).Value2 = '  -13.34%  '

$ws.Cells.Item(22,4).Value2 = "'0.759"
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value2 = '  -12.95%  '

$ws.Cells.Item(23,4).Value2 = "'7.72"
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value2 = '  -14.06%  '

$ws.Cells.Item(24,4).Value2 = "'84.27"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value2 = '  -13.18%  '

$ws.Cells.Item(25,4).Value2 = "'13.30"
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value2 = '  -13.95%  '

$ws.Cells.Item(26,4).Value2 = "'3.32"
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).Value2 = '  -10.25%  '

$ws.Cells.Item(27,5).Value2 = '  -0.12%  '

$ws.Cells.Item(28,2).Value2 = 'EthereumClassic'
$ws.Cells.Item(28,3).Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28,4).Value2 = "'29.14"
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value2 = '  -13.23%  '

$ws.Cells.Item(29,2).Value2 = 'ImmutableX'
$ws.Cells.Item(29,3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(29,4).Value2 = "'2.10"
$ws.Cells.Item(29,4).ClearFormats()
$ws.Cells.Item(29,5).Value2 = '  -17.28%  '

$ws.Cells.Item(30,4).Value2 = "'7.39"
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value2 = '  -17.27%  '

$ws.Cells.Item(31,4).Value2 = "'1.15"
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value2 = '  -9.57%  '

$ws.Cells.Item(32,2).Value2 = 'Bittensor'
$ws.Cells.Item(32,3).Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(32,4).Value2 = "'561.44"
$ws.Cells.Item(32,4).ClearFormats()
$ws.Cells.Item(32,5).Value2 = '  -12.88%  '

$ws.Cells.Item(33,4).Value2 = "'2.46"
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value2 = '  -16.68%  '

$ws.Cells.Item(34,2).Value2 = 'Filecoin'
$ws.Cells.Item(34,3).Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34,4).Value2 = "'6.53"
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value2 = '  -18.26%  '

$ws.Cells.Item(35,4).Value2 = "'5.66"
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value2 = '  -16.87%  '

$ws.Cells.Item(36,4).Value2 = "'1.01"
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value2 = '  +0.80%  '

$ws.Cells.Item(37,4).Value2 = "'54.41"
$ws.Cells.Item(37,4).ClearFormats()
$ws.Cells.Item(37,5).Value2 = '  -3.70%  '

$ws.Cells.Item(38,4).Value2 = "'0.0431"
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value2 = '  -8.01%  '

$ws.Cells.Item(39,2).Value2 = 'Cosmos'
$ws.Cells.Item(39,3).Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(39,4).Value2 = "'9.22"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value2 = '  -13.72%  '

$ws.Cells.Item(40,2).Value2 = 'Hedera'
$ws.Cells.Item(40,3).Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40,4).Value2 = "'0.0843"
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value2 = '  -14.44%  '

$ws.Cells.Item(41,4).Value2 = "'0.128"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value2 = '  -9.24%  '

$ws.Cells.Item(42,2).Value2 = 'Maker'
$ws.Cells.Item(42,3).Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42,4).Value2 = "'2.901.46"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value2 = '  -13.36%  '

$ws.Cells.Item(43,2).Value2 = 'dogwifhat'
$ws.Cells.Item(43,3).Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43,4).Value2 = "'2.64"
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value2 = '  -25.28%  '

$ws.Cells.Item(44,4).Value2 = "'2.45"
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value2 = '  -14.65%  '

$ws.Cells.Item(45,4).Value2 = "'0.257"
$ws.Cells.Item(45,4).ClearFormats()
$ws.Cells.Item(45,5).Value2 = '  -15.97%  '

$ws.Cells.Item(46,2).Value2 = 'USDe'
$ws.Cells.Item(46,3).Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46,4).Value2 = "'1.00"
$ws.Cells.Item(46,4).ClearFormats()
$ws.Cells.Item(46,5).Value2 = '  -0.01%  '

$ws.Cells.Item(47,2).Value2 = 'PEPE'
$ws.Cells.Item(47,3).Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(47,4).Value2 = "'0.0₃0548"
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value2 = '  -25.46%  '

$ws.Cells.Item(48,2).Value2 = 'Monero'
$ws.Cells.Item(48,3).Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(48,4).Value2 = "'125.00"
$ws.Cells.Item(48,4).ClearFormats()
$ws.Cells.Item(48,5).Value2 = '  -6.50%  '

$ws.Cells.Item(49,4).Value2 = "'0.113"
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value2 = '  -11.97%  '

$ws.Cells.Item(50,2).Value2 = 'Fetch.AI'
$ws.Cells.Item(50,3).Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(50,4).Value2 = "'2.06"
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value2 = '  -18.48%  '

$ws.Cells.Item(51,2).Value2 = 'InjectiveProtocol'
$ws.Cells.Item(51,3).Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51,4).Value2 = "'25.16"
$ws.Cells.Item(51,4).ClearFormats()
$ws.Cells.Item(51,5).Value2 = '  -21.46%  '
